# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-record detail columns (D, L, M, N, O,
# P, Q, R, S, T) across the existing data rows (2-23, 26); rows 24-25 keep
# their original data. We snapshot every source row's values first (since the
# remap is a permutation, not a simple append), then write them back out so
# no value is clobbered before it's been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# target row -> source row (values currently sitting in $source row move to $target row)
$map = @{
    2  = 11
    3  = 4
    4  = 10
    5  = 8
    6  = 16
    7  = 2
    8  = 19
    9  = 20
    10 = 14
    11 = 3
    12 = 6
    13 = 7
    14 = 21
    15 = 22
    16 = 9
    17 = 26
    18 = 13
    19 = 15
    20 = 17
    21 = 5
    22 = 23
    23 = 18
    26 = 12
}

# 1) Snapshot current values for every row referenced as a source.
$snapshot = @{}
foreach ($targetRow in $map.Keys) {
    $srcRow = $map[$targetRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$srcRow").Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# 2) Write the snapshotted source-row values into each target row.
foreach ($targetRow in $map.Keys) {
    $srcRow = $map[$targetRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $rowVals[$col]
    }
}
